# Update to "Colombia Primera B" worksheet reflecting a league-data refresh.
#
# The underlying change is a reordering of several fixture rows after a
# data refresh (new rows were inserted / re-sorted upstream), which moved
# the match data that used to live in certain rows into different rows:
#   - Row 10  <->  Row 11
#   - Row 63  <->  Row 64
#   - Row 104 -> Row 105 -> Row 106 -> Row 104   (3-way rotation)
#   - Row 263 <->  Row 264
# The "id" column (A) keeps its original sequential value per row; all
# other columns (match id, teams, scores and odds) move together with the
# match they describe, since a match's data belongs together.
#
# (The workbook's raw sharedStrings table also re-numbers the "Llaneros"
# and "Real Cartagena" entries, but every single cell that referenced
# those two entries keeps showing exactly the same team name afterwards,
# so from the worksheet's point of view nothing else changes - it is purely
# an artifact of how the shared string table got regenerated. Writing the
# cell values below via COM lets Excel pick the right shared-string
# automatically, so we don't need to touch those other rows at all.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colombia Primera B")

    # Row 10
    $ws.Range("B10").Value = 6695291
    $ws.Range("E10").Value = "Real Cartagena"
    $ws.Range("F10").Value = "Fortaleza"
    $ws.Range("G10").Value = 0
    $ws.Range("H10").Value = 1
    $ws.Range("I10").Value = "A"
    $ws.Range("J10").Value = 2.875
    $ws.Range("K10").Value = 2.9
    $ws.Range("L10").Value = 2.375
    $ws.Range("M10").Value = 3.8
    $ws.Range("N10").Value = 3.1
    $ws.Range("O10").Value = 1.909
    $ws.Range("P10").Value = 0.5
    $ws.Range("Q10").Value = 1.825
    $ws.Range("R10").Value = 1.975
    $ws.Range("S10").Value = 2.25
    $ws.Range("T10").Value = 2
    $ws.Range("U10").Value = 1.8
    $ws.Range("V10").Value = -1
    $ws.Range("W10").Value = -1
    $ws.Range("X10").Value = 0.909
    $ws.Range("Y10").Value = -1
    $ws.Range("Z10").Value = 0.9750000000000001
    $ws.Range("AA10").Value = -1
    $ws.Range("AB10").Value = 0.8

    # Row 11
    $ws.Range("B11").Value = 6695290
    $ws.Range("E11").Value = "Llaneros"
    $ws.Range("F11").Value = "Deportes Quindio"
    $ws.Range("G11").Value = 2
    $ws.Range("H11").Value = 0
    $ws.Range("I11").Value = "H"
    $ws.Range("J11").Value = 1.666
    $ws.Range("K11").Value = 3.2
    $ws.Range("L11").Value = 5
    $ws.Range("M11").Value = 1.45
    $ws.Range("N11").Value = 3.6
    $ws.Range("O11").Value = 6.5
    $ws.Range("P11").Value = -1
    $ws.Range("Q11").Value = 1.825
    $ws.Range("R11").Value = 1.975
    $ws.Range("S11").Value = 2.25
    $ws.Range("T11").Value = 1.925
    $ws.Range("U11").Value = 1.875
    $ws.Range("V11").Value = 0.45
    $ws.Range("W11").Value = -1
    $ws.Range("X11").Value = -1
    $ws.Range("Y11").Value = 0.825
    $ws.Range("Z11").Value = -1
    $ws.Range("AA11").Value = -0.5
    $ws.Range("AB11").Value = 0.4375

    # Row 63
    $ws.Range("B63").Value = 6990754
    $ws.Range("E63").Value = "Tigres FC"
    $ws.Range("F63").Value = "Llaneros"
    $ws.Range("G63").Value = 0
    $ws.Range("H63").Value = 1
    $ws.Range("I63").Value = "A"
    $ws.Range("J63").Value = 2.9
    $ws.Range("K63").Value = 3.2
    $ws.Range("L63").Value = 2.2
    $ws.Range("M63").Value = 2.4
    $ws.Range("N63").Value = 3.25
    $ws.Range("O63").Value = 3
    $ws.Range("P63").Value = 0
    $ws.Range("Q63").Value = 1.7
    $ws.Range("R63").Value = 2
    $ws.Range("S63").Value = 2
    $ws.Range("T63").Value = 2.025
    $ws.Range("U63").Value = 1.775
    $ws.Range("V63").Value = -1
    $ws.Range("W63").Value = -1
    $ws.Range("X63").Value = 2
    $ws.Range("Y63").Value = -1
    $ws.Range("Z63").Value = 1
    $ws.Range("AA63").Value = -1
    $ws.Range("AB63").Value = 0.7749999999999999

    # Row 64
    $ws.Range("B64").Value = 7121333
    $ws.Range("E64").Value = "Atletico Cali FC"
    $ws.Range("F64").Value = "Real Soacha Cundinamarca"
    $ws.Range("G64").Value = 1
    $ws.Range("H64").Value = 2
    $ws.Range("I64").Value = "A"
    $ws.Range("J64").Value = 2.25
    $ws.Range("K64").Value = 3
    $ws.Range("L64").Value = 3
    $ws.Range("M64").Value = 2.375
    $ws.Range("N64").Value = 3.2
    $ws.Range("O64").Value = 3.2
    $ws.Range("P64").Value = -0.25
    $ws.Range("Q64").Value = 2
    $ws.Range("R64").Value = 1.8
    $ws.Range("S64").Value = 2
    $ws.Range("T64").Value = 1.8
    $ws.Range("U64").Value = 2
    $ws.Range("V64").Value = -1
    $ws.Range("W64").Value = -1
    $ws.Range("X64").Value = 2.2
    $ws.Range("Y64").Value = -1
    $ws.Range("Z64").Value = 0.8
    $ws.Range("AA64").Value = 0.8
    $ws.Range("AB64").Value = -1

    # Row 104
    $ws.Range("B104").Value = 6990794
    $ws.Range("E104").Value = "Real San Andres"
    $ws.Range("F104").Value = "Tigres FC"
    $ws.Range("G104").Value = 1
    $ws.Range("H104").Value = 2
    $ws.Range("I104").Value = "A"
    $ws.Range("J104").Value = 1.909
    $ws.Range("K104").Value = 3.2
    $ws.Range("L104").Value = 3.75
    $ws.Range("M104").Value = 1.75
    $ws.Range("N104").Value = 3.5
    $ws.Range("O104").Value = 5
    $ws.Range("P104").Value = -0.5
    $ws.Range("Q104").Value = 1.75
    $ws.Range("R104").Value = 2.05
    $ws.Range("S104").Value = 2.25
    $ws.Range("T104").Value = 1.85
    $ws.Range("U104").Value = 1.95
    $ws.Range("V104").Value = -1
    $ws.Range("W104").Value = -1
    $ws.Range("X104").Value = 4
    $ws.Range("Y104").Value = -1
    $ws.Range("Z104").Value = 1.05
    $ws.Range("AA104").Value = 0.8500000000000001
    $ws.Range("AB104").Value = -1

    # Row 105
    $ws.Range("B105").Value = 6990792
    $ws.Range("E105").Value = "Fortaleza"
    $ws.Range("F105").Value = "Real Cartagena"
    $ws.Range("G105").Value = 1
    $ws.Range("H105").Value = 0
    $ws.Range("I105").Value = "H"
    $ws.Range("J105").Value = 1.8
    $ws.Range("K105").Value = 3.3
    $ws.Range("L105").Value = 4
    $ws.Range("M105").Value = 1.571
    $ws.Range("N105").Value = 4
    $ws.Range("O105").Value = 6.5
    $ws.Range("P105").Value = -1
    $ws.Range("Q105").Value = 1.95
    $ws.Range("R105").Value = 1.85
    $ws.Range("S105").Value = 2.25
    $ws.Range("T105").Value = 1.85
    $ws.Range("U105").Value = 1.95
    $ws.Range("V105").Value = 0.571
    $ws.Range("W105").Value = -1
    $ws.Range("X105").Value = -1
    $ws.Range("Y105").Value = 0
    $ws.Range("Z105").Value = 0
    $ws.Range("AA105").Value = -1
    $ws.Range("AB105").Value = 0.95

    # Row 106
    $ws.Range("B106").Value = 6990788
    $ws.Range("E106").Value = "Boca Juniors De Cali"
    $ws.Range("F106").Value = "Cucuta Deportivo"
    $ws.Range("G106").Value = 1
    $ws.Range("H106").Value = 0
    $ws.Range("I106").Value = "H"
    $ws.Range("J106").Value = 3
    $ws.Range("K106").Value = 3
    $ws.Range("L106").Value = 2.3
    $ws.Range("M106").Value = 3
    $ws.Range("N106").Value = 3
    $ws.Range("O106").Value = 2.5
    $ws.Range("P106").Value = 0
    $ws.Range("Q106").Value = 2.05
    $ws.Range("R106").Value = 1.75
    $ws.Range("S106").Value = 2
    $ws.Range("T106").Value = 1.85
    $ws.Range("U106").Value = 1.95
    $ws.Range("V106").Value = 2
    $ws.Range("W106").Value = -1
    $ws.Range("X106").Value = -1
    $ws.Range("Y106").Value = 1.05
    $ws.Range("Z106").Value = -1
    $ws.Range("AA106").Value = -1
    $ws.Range("AB106").Value = 0.95

    # Row 263
    $ws.Range("B263").Value = 7657928
    $ws.Range("E263").Value = "Atletico Huila"
    $ws.Range("F263").Value = "Barranquilla FC"
    $ws.Range("G263").Value = 2
    $ws.Range("H263").Value = 0
    $ws.Range("I263").Value = "H"
    $ws.Range("J263").Value = 1.444
    $ws.Range("K263").Value = 4.333
    $ws.Range("L263").Value = 6
    $ws.Range("M263").Value = 1.571
    $ws.Range("N263").Value = 4.2
    $ws.Range("O263").Value = 5.25
    $ws.Range("P263").Value = -1
    $ws.Range("Q263").Value = 1.975
    $ws.Range("R263").Value = 1.825
    $ws.Range("S263").Value = 2.5
    $ws.Range("T263").Value = 1.925
    $ws.Range("U263").Value = 1.875
    $ws.Range("V263").Value = 0.571
    $ws.Range("W263").Value = -1
    $ws.Range("X263").Value = -1
    $ws.Range("Y263").Value = 0.9750000000000001
    $ws.Range("Z263").Value = -1
    $ws.Range("AA263").Value = -1
    $ws.Range("AB263").Value = 0.875

    # Row 264
    $ws.Range("B264").Value = 7658173
    $ws.Range("E264").Value = "Atletico Cali FC"
    $ws.Range("F264").Value = "Deportes Quindio"
    $ws.Range("G264").Value = 0
    $ws.Range("H264").Value = 1
    $ws.Range("I264").Value = "A"
    $ws.Range("J264").Value = 4.5
    $ws.Range("K264").Value = 3.5
    $ws.Range("L264").Value = 1.727
    $ws.Range("M264").Value = 7.5
    $ws.Range("N264").Value = 5.5
    $ws.Range("O264").Value = 1.363
    $ws.Range("P264").Value = 1.25
    $ws.Range("Q264").Value = 1.9
    $ws.Range("R264").Value = 1.9
    $ws.Range("S264").Value = 3
    $ws.Range("T264").Value = 1.9
    $ws.Range("U264").Value = 1.9
    $ws.Range("V264").Value = -1
    $ws.Range("W264").Value = -1
    $ws.Range("X264").Value = 0.363
    $ws.Range("Y264").Value = 0.45
    $ws.Range("Z264").Value = -0.5
    $ws.Range("AA264").Value = -1
    $ws.Range("AB264").Value = 0.8999999999999999

